$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("19Tto20TMap")

# Populate the new "pe" column (M) for all data rows (2-77) with "2018Oct"
$ws.Range("M2:M77").Value = "2018Oct"

# Reset the view: scroll back to A1 and select A1 (matches author re-saving
# the file after scrolling back / clearing the stray selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select() | Out-Null
